# Regression Project presentation.pptx - apply "Link to Project" link restyle
#
# Slide 1, shape "Subtitle 2" has a paragraph that reads "Link to Project"
# made up of two hyperlinked runs:
#   "Link to " -> rId2 (kaggle link)
#   "Project"  -> rId3 (github link)
#
# The edit re-splits that text into three runs with new colors while keeping
# the same overall wording, and re-points the first two runs at rId2:
#   "Link "  -> rId2, color F49100
#   "to "    -> rId2, color FF9900
#   "Project"-> rId3, color FF9900

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

$kaggleUrl = "https://www.kaggle.com/ronishternberg/roni-california-housing-prices"
$githubUrl = "https://github.com/linoybar/Regression_Python_Project"

# Locate the "Link to Project" paragraph (5th paragraph in the text box) and
# work out the character offsets for the three new runs.
$para = $tr.Paragraphs(5)
$paraStart = $para.Start

$runLinkStart = $paraStart
$runToStart = $paraStart + 5
$runProjectStart = $paraStart + 8

$runLink = $tr.Characters($runLinkStart, 5)
$runTo = $tr.Characters($runToStart, 3)
$runProject = $tr.Characters($runProjectStart, 7)

# "Link " - orange (F49100), still linking to the kaggle notebook (rId2)
$runLink.Font.Color.RGB = 37364
$linkAction = $runLink.ActionSettings.Item(1)
$linkAction.Hyperlink.Address = $kaggleUrl

# "to " - now also linking to the kaggle notebook (rId2), brighter orange (FF9900)
$runTo.Font.Color.RGB = 39423
$toAction = $runTo.ActionSettings.Item(1)
$toAction.Hyperlink.Address = $kaggleUrl

# "Project" - keeps linking to the github repo (rId3), brighter orange (FF9900)
$runProject.Font.Color.RGB = 39423
$projectAction = $runProject.ActionSettings.Item(1)
$projectAction.Hyperlink.Address = $githubUrl
